$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    479, 745, 1022, 1576, 2610, 2859, 3552, 4149, 4266, 4436,
    4553, 5827, 6296, 7470, 8163, 10071, 10206, 10483, 11079, 11836,
    12273, 14503, 14940, 15280, 15685, 16220, 16497, 16860, 17382, 17995,
    18272, 19572, 20393, 20606, 22226, 23827, 24125, 24999, 25617, 28257,
    28262, 28411, 29445, 31886, 32346, 32911, 34030, 34467, 35394, 39856,
    40570, 40751, 42236, 42993, 46147, 48262, 49146, 51183, 52664, 53223,
    53713, 54306, 56075, 56160, 58996, 59209, 61075, 61565, 63052, 67374,
    68792, 69815, 75305, 76825, 80106, 82822, 84997, 88502, 89227, 91190,
    95169, 96259, 96504, 97698, 98692, 102594
)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $newValues[$i]
}
